$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 7
$ws.Range("B7").Value = 0.6743055555555556
$ws.Range("C7").Value = 5321

# Update value in row 12
$ws.Range("C12").Value = 0.67569444444444438

# Update the selected cell to match the diff
$ws.Range("G13").Select()
